# Generate Report for Archive
# The localization status for "eaf4563a-563d-4ade-8eb6-56a9aa18b524.md" moved
# forward in its workflow: it is now back "In Translation" (no longer
# "Ready for handoff") for both the zh-cn and de-de targets. Update the
# per-language status tables and the Overview summary sheet accordingly.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# zh-cn sheet: row 5 (header row 1, data rows 2-5) is the
# eaf4563a-563d-4ade-8eb6-56a9aa18b524.md entry -> Status column C
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C4").Value = $newStatus

# de-de sheet: same row/column layout
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C4").Value = $newStatus

# Overview sheet: columns E (zh-cn) and F (de-de) on the same file's row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus
